# The commit removes the last slide (sldId="272", r:id="rId7", the slide
# backed by ppt/slides/slide6.xml) from the deck. Deleting it through the
# Slides collection removes the slide part, its relationship, the
# <p:sldId> entry in the presentation's slide list, and the
# [Content_Types].xml override, mirroring the target diff.
$p = $ppt.ActivePresentation
$p.Slides.Item(6).Delete()
